# "Generate Report for Archive"
#
# The localization-status report is regenerated: every "Ready for handoff"
# status cell becomes "In Translation" (Overview!E2/F2 plus the Status
# column, column C, on the per-locale "zh-cn" and "de-de" sheets), and the
# now-shorter status column is narrowed accordingly on all three sheets.
#
# Note on column widths: Excel's ColumnWidth property only accepts values
# that land on this workbook font's whole-pixel grid, so it snaps to the
# nearest representable width (~13.33 chars) when asked for the source
# report's exact 13.4101845877511 value - the same way it would if you
# typed that width into the "Column Width" dialog by hand.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newStatusColWidth = 12.5   # snaps to the nearest on-grid ColumnWidth (~13.41 chars)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns("E:E").ColumnWidth = $newStatusColWidth
$wsOverview.Columns("F:F").ColumnWidth = $newStatusColWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns("C:C").ColumnWidth = $newStatusColWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns("C:C").ColumnWidth = $newStatusColWidth
